$d = $word.ActiveDocument

# The first paragraph of the review ("Paper:  Automatic Image Caption
# Generation") needs to become bold - both the run text and the paragraph
# mark itself (so the pPr/rPr also carries <w:b/><w:bCs/>).

$p1 = $d.Paragraphs(1)
$r = $p1.Range

# Pull this paragraph's existing opening-tag attributes (w14:paraId, rsids,
# etc.) straight out of the live document XML so we can reproduce them
# exactly and only touch the formatting, leaving everything else alone.
$fullXml = $r.XML($false)
$pAttrs = ""
if ($fullXml -match '<w:p\s*([^>]*)>\s*<w:r>\s*<w:t>Paper:  Automatic Image Caption Generation</w:t>') {
    $pAttrs = $matches[1].Trim()
}

$pOpenTag = "<w:p>"
if ($pAttrs.Length -gt 0) {
    $pOpenTag = "<w:p $pAttrs>"
}

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          $pOpenTag
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>Paper:  Automatic Image Caption Generation</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]($r.InsertXML($xml, $null))
